$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 1.34
$ws.Range("P2").Value = 2.16
$ws.Range("S2").Value = 3

# Row 3
$ws.Range("J3").Value = 3.55
$ws.Range("N3").Value = 2

# Row 4
$ws.Range("T4").Value = 1.67
$ws.Range("U4").Value = 2.2
$ws.Range("Z4").Value = 90

# Row 6
$ws.Range("W6").Value = 2.06

# Row 7
$ws.Range("I7").Value = 6.2

# Row 8
$ws.Range("P8").Value = 1.73
$ws.Range("T8").Value = 1.82
$ws.Range("U8").Value = 1.86
$ws.Range("AF8").Value = 21
$ws.Range("AG8").Value = 17
$ws.Range("AL8").Value = 65

# Row 9
$ws.Range("G9").Value = 2.14
$ws.Range("N9").Value = 1.72
$ws.Range("P9").Value = 1.72
$ws.Range("Q9").Value = 2.06
$ws.Range("W9").Value = 1.87

$wb.Save()
